# Auto-generated edit script applying the diff changes to Maduin_Profits workbook
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 559.38464
$ws.Range("I53").Value = 176
$ws.Range("K53").Value = 176
$ws.Range("M53").Value = 461
# Row 107
$ws.Range("H107").Value = 616.75
$ws.Range("I107").Value = 641.2
$ws.Range("J107").Value = 250
$ws.Range("K107").Value = 641.2
$ws.Range("L107").Value = 250
$ws.Range("M107").Value = 1278.8
$ws.Range("N107").Value = -4090
# Row 125
$ws.Range("H125").Value = 4889
$ws.Range("I125").Value = 2732
$ws.Range("J125").Value = 5967.5
$ws.Range("K125").Value = 24588
$ws.Range("L125").Value = 53707.5
$ws.Range("M125").Value = -22128
$ws.Range("N125").Value = -58627.5
# Row 138
$ws.Range("H138").Value = 2702.2856
$ws.Range("I138").Value = 2702.2856
$ws.Range("K138").Value = 8106.8568
$ws.Range("M138").Value = -2966.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 4374.5
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 5666
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 5666
$ws.Range("M12").Value = -327
$ws.Range("N12").Value = -6012
# Row 28
$ws.Range("H28").Value = 14133.4
$ws.Range("I28").Value = 14133.4
$ws.Range("K28").Value = 14133.4
$ws.Range("M28").Value = -13941.4
# Row 74
$ws.Range("H74").Value = 935.35
$ws.Range("I74").Value = 935.35
$ws.Range("K74").Value = 935.35
$ws.Range("M74").Value = -61.35000000000002
# Row 77
$ws.Range("H77").Value = 935.35
$ws.Range("I77").Value = 935.35
$ws.Range("K77").Value = 4676.75
$ws.Range("M77").Value = -308.75
# Row 99
$ws.Range("H99").Value = 14133.4
$ws.Range("I99").Value = 14133.4
$ws.Range("K99").Value = 14133.4
$ws.Range("M99").Value = -11138.4
# Row 122
$ws.Range("H122").Value = 3084
$ws.Range("J122").Value = 4492.3335
$ws.Range("L122").Value = 13477.0005
$ws.Range("N122").Value = -18377.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 8040
$ws.Range("I20").Value = 6400
$ws.Range("J20").Value = 10500
$ws.Range("K20").Value = 6400
$ws.Range("L20").Value = 10500
$ws.Range("M20").Value = -6153
$ws.Range("N20").Value = -10994
# Row 86
$ws.Range("H86").Value = 4832.6665
$ws.Range("J86").Value = 4197.5713
$ws.Range("L86").Value = 4197.5713
$ws.Range("N86").Value = -6443.5713
# Row 89
$ws.Range("H89").Value = 4832.6665
$ws.Range("J89").Value = 4197.5713
$ws.Range("L89").Value = 20987.8565
$ws.Range("N89").Value = -32219.8565
# Row 99
$ws.Range("H99").Value = 3185.6924
$ws.Range("I99").Value = 2113
$ws.Range("J99").Value = 5599.25
$ws.Range("K99").Value = 2113
$ws.Range("L99").Value = 5599.25
$ws.Range("M99").Value = -615
$ws.Range("N99").Value = -8595.25
# Row 105
$ws.Range("H105").Value = 2435.7896
$ws.Range("I105").Value = 2397.4
$ws.Range("J105").Value = 2579.75
$ws.Range("K105").Value = 2397.4
$ws.Range("L105").Value = 2579.75
$ws.Range("M105").Value = -650.4000000000001
$ws.Range("N105").Value = -6073.75
# Row 134
$ws.Range("H134").Value = 1820.9474
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 2377.4666
$ws.Range("I7").Value = 1205.6364
$ws.Range("J7").Value = 5600
$ws.Range("K7").Value = 1205.6364
$ws.Range("L7").Value = 5600
$ws.Range("M7").Value = -1092.6364
$ws.Range("N7").Value = -5826
# Row 53
$ws.Range("H53").Value = 9342
$ws.Range("J53").Value = 9342
$ws.Range("L53").Value = 9342
$ws.Range("N53").Value = -10556
# Row 58
$ws.Range("H58").Value = 1218.9445
$ws.Range("I58").Value = 781.5714
$ws.Range("K58").Value = 781.5714
$ws.Range("M58").Value = -578.5714
# Row 62
$ws.Range("H62").Value = 2250
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
# Row 65
$ws.Range("H65").Value = 2250
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
# Row 132
$ws.Range("H132").Value = 3944.9375
$ws.Range("I132").Value = 4014.6
$ws.Range("K132").Value = 12043.8
$ws.Range("M132").Value = -9513.799999999999
# Row 134
$ws.Range("H134").Value = 1697.125
$ws.Range("I134").Value = 1524.2858
$ws.Range("K134").Value = 4572.857400000001
$ws.Range("M134").Value = -2037.857400000001
# Row 136
$ws.Range("H136").Value = 1218.9445
$ws.Range("I136").Value = 781.5714
$ws.Range("K136").Value = 2344.7142
$ws.Range("M136").Value = 205.2857999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 99.333336
$ws.Range("I2").Value = 54.63158
$ws.Range("J2").Value = 176.54546
$ws.Range("K2").Value = 327.78948
$ws.Range("L2").Value = 1059.27276
$ws.Range("M2").Value = -214.78948
$ws.Range("N2").Value = -1285.27276
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0
# Row 68
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
# Row 71
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10299
$ws.Range("I70").Value = 10498.6
$ws.Range("K70").Value = 10498.6
$ws.Range("M70").Value = -10228.6
# Row 73
$ws.Range("H73").Value = 10299
$ws.Range("I73").Value = 10498.6
$ws.Range("K73").Value = 10498.6
$ws.Range("M73").Value = -9562.6
# Row 97
$ws.Range("H97").Value = 953.2222
$ws.Range("I97").Value = 509.875
$ws.Range("K97").Value = 509.875
$ws.Range("M97").Value = -13.875
# Row 132
$ws.Range("H132").Value = 3198.818
$ws.Range("I132").Value = 3328.4285
$ws.Range("J132").Value = 2972
$ws.Range("K132").Value = 9985.2855
$ws.Range("L132").Value = 8916
$ws.Range("M132").Value = -7455.2855
$ws.Range("N132").Value = -13976

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 379.6
$ws.Range("J22").Value = 299.33334
$ws.Range("L22").Value = 299.33334
$ws.Range("N22").Value = -889.33334
# Row 27
$ws.Range("H27").Value = 379.6
$ws.Range("J27").Value = 299.33334
$ws.Range("L27").Value = 299.33334
$ws.Range("N27").Value = -513.33334
# Row 46
$ws.Range("H46").Value = 3866.2
$ws.Range("J46").Value = 3704.5454
$ws.Range("L46").Value = 3704.5454
$ws.Range("N46").Value = -4080.5454
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 106
$ws.Range("H106").Value = 11683.25
$ws.Range("J106").Value = 11683.25
$ws.Range("L106").Value = 11683.25
$ws.Range("N106").Value = -14207.25

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 23999.666
$ws.Range("J18").Value = 23999.666
$ws.Range("L18").Value = 23999.666
$ws.Range("N18").Value = -24345.666
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0
# Row 69
$ws.Range("H69").Value = 108544.2
$ws.Range("J69").Value = 108544.2
$ws.Range("L69").Value = 108544.2
$ws.Range("N69").Value = -110042.2
# Row 72
$ws.Range("H72").Value = 108544.2
$ws.Range("J72").Value = 108544.2
$ws.Range("L72").Value = 325632.6
$ws.Range("N72").Value = -333120.6
# Row 100
$ws.Range("H100").Value = 4979703
$ws.Range("I100").Value = 6970315
$ws.Range("J100").Value = 3173.75
$ws.Range("K100").Value = 13940630
$ws.Range("L100").Value = 6347.5
$ws.Range("M100").Value = -13940089
$ws.Range("N100").Value = -7429.5
# Row 136
$ws.Range("H136").Value = 1499.3334
$ws.Range("I136").Value = 1499.3334
$ws.Range("K136").Value = 4498.0002
$ws.Range("M136").Value = -1948.0002
